$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free inline approach: force each touched cell to Text format before
# assigning its literal string value (so values like "29.433.20", "66.13",
# "0.00000000119" are preserved exactly as text, matching the original
# inlineStr cell type), then restore the default "Normal" style so no stray
# number-format styling is left behind on the cell.

function Set-TextValue($rangeAddr, $val) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '29.433.20'
Set-TextValue 'D3' '1.855.59'
Set-TextValue 'E3' '  +0.39%  '
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '241.53'
Set-TextValue 'E5' '  +0.28%  '
Set-TextValue 'D6' '0.6329'
Set-TextValue 'E6' '  +0.82%  '
Set-TextValue 'D8' '0.07591'
Set-TextValue 'D9' '0.2927'
Set-TextValue 'E9' '  +0.24%  '
Set-TextValue 'D10' '24.62'
Set-TextValue 'E10' '  -0.84%  '
Set-TextValue 'D11' '0.07761'
Set-TextValue 'E11' '  +0.19%  '
Set-TextValue 'D12' '1.855.79'
Set-TextValue 'E12' '  +0.48%  '
Set-TextValue 'D13' '5.038'
Set-TextValue 'E13' '  +0.14%  '
Set-TextValue 'D14' '0.6857'
Set-TextValue 'E14' '  +0.73%  '
Set-TextValue 'D15' '0.00001045'
Set-TextValue 'E15' '  -2.52%  '
Set-TextValue 'D16' '83.37'
Set-TextValue 'E16' '  -0.04%  '
Set-TextValue 'D17' '2.115.02'
Set-TextValue 'E17' '  +0.70%  '
Set-TextValue 'D18' '6.152'
Set-TextValue 'E18' '  -0.37%  '
Set-TextValue 'D19' '29.455.08'
Set-TextValue 'E19' '  +0.02%  '
Set-TextValue 'D20' '230.52'
Set-TextValue 'E20' '  +0.98%  '
Set-TextValue 'E21' '  -0.06%  '
Set-TextValue 'E22' '  +0.01%  '
Set-TextValue 'D23' '7.530'
Set-TextValue 'E23' '  +1.48%  '
Set-TextValue 'E24' '  +0.01%  '
Set-TextValue 'D25' '159.18'
Set-TextValue 'E25' '  +0.77%  '
Set-TextValue 'D26' '0.1399'
Set-TextValue 'E26' '  +1.73%  '
Set-TextValue 'D27' '8.483'
Set-TextValue 'E27' '  +0.90%  '
Set-TextValue 'D28' '17.75'
Set-TextValue 'E28' '  +0.37%  '
Set-TextValue 'D29' '1.417'
Set-TextValue 'E29' '  +5.09%  '
Set-TextValue 'D30' '1.484'
Set-TextValue 'E30' '  +1.26%  '
Set-TextValue 'D31' '0.05686'
Set-TextValue 'E31' '  +0.30%  '
Set-TextValue 'D32' '4.158'
Set-TextValue 'E32' '  +0.86%  '
Set-TextValue 'D33' '4.064'
Set-TextValue 'E33' '  +0.93%  '
Set-TextValue 'D34' '1.832'
Set-TextValue 'E34' '  -0.60%  '
Set-TextValue 'D35' '1.158'
Set-TextValue 'E35' '  -0.38%  '
Set-TextValue 'D36' '0.6985'
Set-TextValue 'E36' '  -0.53%  '
Set-TextValue 'D37' '2.589'
Set-TextValue 'E37' '  -0.01%  '
Set-TextValue 'D38' '1.251.71'
Set-TextValue 'E38' '  +2.05%  '
Set-TextValue 'E39' '  +2.24%  '
Set-TextValue 'D40' '2.775'
Set-TextValue 'E40' '  +0.32%  '
Set-TextValue 'D41' '6.520'
Set-TextValue 'E41' '  -0.36%  '
Set-TextValue 'D42' '0.9110'
Set-TextValue 'E42' '  +0.72%  '
Set-TextValue 'E43' '  +0.02%  '
Set-TextValue 'D44' '2.018.45'
Set-TextValue 'E44' '  +0.49%  '
Set-TextValue 'E45' '  -0.23%  '
Set-TextValue 'D46' '66.13'
Set-TextValue 'E46' '  +0.12%  '
Set-TextValue 'D47' '7.165'
Set-TextValue 'B48' 'BabyDogeCoin'
Set-TextValue 'C48' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D48' '0.00000000119'
Set-TextValue 'E48' '  -2.12%  '
Set-TextValue 'B49' 'Algorand'
Set-TextValue 'C49' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D49' '0.1167'
Set-TextValue 'E49' '  +1.05%  '
Set-TextValue 'B50' 'EnergySwap'
Set-TextValue 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '9.068'
Set-TextValue 'E50' '  +0.61%  '
Set-TextValue 'B51' 'TheSandbox'
Set-TextValue 'C51' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D51' '0.3970'
Set-TextValue 'E51' '  -1.18%  '
